$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Work from the bottom of the sheet upward so that earlier (lower) row
# numbers used below are not disturbed by inserts/deletes done later.
# ---------------------------------------------------------------------

# 1) Insert six new rows after the "i-racer" row (row 15) for the new
#    Arduino / motors / LEDs / resistors / controllers / keyboard items.
$ws.Rows("16:21").Insert()

$ws.Range("A16").Value = "Arduino Mega 2560"
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = "Arduino"

$ws.Range("A17").Value = "DC Motors"
$ws.Range("B17").Value = 4
$ws.Range("C17").Value = "Various motors"

$ws.Range("A18").Value = "LEDs "
$ws.Range("B18").Value = "x"
$ws.Range("C18").Value = "Blue, Yellow, Green and Red"

$ws.Range("A19").Value = "Resistors"
$ws.Range("B19").Value = "x"
$ws.Range("C19").Value = "Varried resistors"

$ws.Range("A20").Value = "Game controllers"
$ws.Range("B20").Value = 4
$ws.Range("C20").Value = "Plus 2 that don't work well"

$ws.Range("A21").Value = "Wireless Keyboard/Air mouse"
$ws.Range("B21").Value = 2
$ws.Range("C21").Value = "Keyboard and mouse in the form factor of a remote control"
$ws.Rows("21").RowHeight = 26

# 2) Replace the "RaspiRobot" row (row 10) with three new rows: battery
#    connectors, variable voltage power supply and power brick.
$ws.Rows("10").Delete()
$ws.Rows("10:12").Insert()

$ws.Range("A10").Value = "Various battery connectors"
$ws.Range("C10").Value = "9V, 2 AA, 4 AA"
$ws.Range("B10").Clear()

$ws.Range("A11").Value = "Variable voltage power supply"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "power small external components"

$ws.Range("A12").Value = "Variable voltage power brick"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "12V-24V up to 5A"
$ws.Rows("12").RowHeight = 13

# 3) Insert the new "chipKIT Pi" row before the "4.3\" Display" row (row 5).
$ws.Rows("5").Insert()
$ws.Range("A5").Value = "chipKIT Pi"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "PIC32 microcontroller"

# 4) Bump the Raspberry Pi B count from 4 to 5.
$ws.Range("B2").Value = 5

# ---------------------------------------------------------------------
# Selection, as recorded in the edited workbook.
# ---------------------------------------------------------------------
$ws.Range("C7").Select()
